$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.150.39'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '2.272.15'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '305.72'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '93.55'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.38%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.489'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('B10').Value = 'BinanceUSD'
$ws.Range('C10').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '30.02'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2,899.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '32.96'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0803'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '2.623.59'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.33'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('D17').Value = '2.275.90'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.786'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.86%  '
$ws.Range('D19').Value = '42.016.34'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.69'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.36%  '
$ws.Range('D21').Value = '0.0₃0920'
$ws.Range('E21').Value = '  +1.98%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.99'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.23'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '244.33'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.61'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '24.03'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.70'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('E30').Value = '  -9.02%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '35.30'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.24%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '160.40'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('E33').Value = '  +3.34%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0744'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.05'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '17.20'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +4.18%  '
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.106'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.04'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.97%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.014.58'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '19.67'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('E45').Value = '  +9.77%  '
$ws.Range('E46').Value = '  +1.63%  '
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.91'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '53.31'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.45%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '72.77'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('E51').Value = '  +0.36%  '
